$d = $word.ActiveDocument

# --- Heading paragraph: title line + huggingface link line ---
$d.Content.Find.Execute(
    "Review 141: [Short] DOLA: DECODING BY CONTRASTING LAYERS IMPROVES FACTUALITY IN LARGE LANGUAGE MODELS, 08.09.2023",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Review 140: [Short] One Wide Feedforward is All You Need, 07.09.2023", 2) | Out-Null

$d.Content.Find.Execute(
    "https://huggingface.co/papers/2309.03883",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://huggingface.co/papers/2309.01826", 2) | Out-Null

# --- Bold "Paper:" line ---
$d.Content.Find.Execute(
    "Paper: https://arxiv.org/abs/2309.03883v2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Paper: https://arxiv.org/abs/2309.01826v2", 2) | Out-Null

# --- Body paragraph: replace all 5 text runs with the new 4 text runs.        ---
# --- Using `v (vertical tab) as the literal manual-line-break marker so that  ---
# --- Word re-creates the <w:br/> separated <w:t> runs exactly like the diff. ---
$newBody = "ארכיטקטורת הטרנספורמרים היא המלכה הבלתי מעורערת של עולם AI. רוב המודלים כמו מודלי שפה או מודלי דיפוזיה גנרטיביים המככבים היום בחדשות AI בנויים על הארכיטקטורה הזו. כמובן שיש לא מעט מחקר גם באקדמיה וגם בתעשיה על שיפור ביצועי הטרנספורמרים. `v`vאז היום ב-shorthebrewpapereviews נסקור מאמר שמנסה לשפר שני ההיבטים של הטרנספורמרים: נפח האחסון וכמות חישובים (בכיוון הקטנתם). קודם כל ניזכר שכל בלוק של טרנספורמר (שהוא גרעין של כל מודל המבוסס על הטרנספורמרים) בנוי ממנגנון של תשומת הלב (attention) ועוד שתי שכבות של fully-connected שאחת מהן עם ReLU והשנייה לינארית. `v`vלפי המאמר השכבות האלו מהוות 2/3 ממספר המשקלים (ב-BERT) וכמובן ״תורמים״ לעומס החישובי. המחברים שאולים מה יקרה עם נוותר על השכבות האלו או שנעשה אותם ״שיתופיים״ (shared) בין כל בלוקי הטרנספורמים של המודל. זה עתיד להקטין את כמות המשקלים במודל באופן משמעותי כי רוב המודלים מכילים עשרות רבות או מאות בלוקי הטרנספורמרים. `v`vהמאמר גם מציע ״לשתף״ משקלים בין האנקודר לדקודר. אבל איך זה משפיע על ביצועים. המאמר מראה שהפגיעה בביצועים לא גדולה במיוחד (למרות שהם ביצעו מספר בדיקות די מצומצם והם בדקו זאת על מודלים די קטנים עם 6 בלוקי טרנספורמרים בלבד). הם גם השווה דמיון בין הייצוגים של משפטים עבור המודל המקורי והמודל ״הקל״ המוצע וגילו שהוא די גבוה. בנוסף הם השווה k משפטים הדומים ביותר מבחינת הייצוג לשני המודלים וגילו גם כאן דמיון רב. נראה מבטיח אך נדרשות בדיקות מקיפות יותר על מודלים רציניים יותר."

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("אנחנו משתמשים במודלי שפה")) {
        $p.Range.Text = $newBody
        break
    }
}

# --- Remove the now-unwanted trailing empty paragraph ---
$count = $d.Paragraphs.Count
$d.Paragraphs($count).Range.Delete()
